$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "'42.079.07"
$ws.Range("E2").Value = "  -0.57%  "
$ws.Range("D3").Formula = "'2.250.84"
$ws.Range("E3").Value = "  -1.06%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Formula = "'306.41"
$ws.Range("E5").Value = "  -0.20%  "
$ws.Range("D6").Formula = "'96.52"
$ws.Range("E6").Value = "  -1.57%  "
$ws.Range("E7").Value = "  -1.04%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  -1.22%  "
$ws.Range("D10").Formula = "'34.61"
$ws.Range("E10").Value = "  -3.16%  "
$ws.Range("D11").Formula = "'0.0815"
$ws.Range("E11").Value = "  +2.22%  "
$ws.Range("D12").Formula = "'0.112"
$ws.Range("E12").Value = "  +0.36%  "
$ws.Range("D13").Formula = "'6.80"
$ws.Range("E13").Value = "  +1.74%  "
$ws.Range("D14").Formula = "'2.600.97"
$ws.Range("E14").Value = "  -1.02%  "
$ws.Range("D15").Formula = "'14.47"
$ws.Range("E15").Value = "  +0.44%  "
$ws.Range("D16").Formula = "'2.254.65"
$ws.Range("E16").Value = "  -0.43%  "
$ws.Range("D17").Formula = "'0.779"
$ws.Range("E17").Value = "  -1.74%  "
$ws.Range("D18").Formula = "'41.949.48"
$ws.Range("E18").Value = "  -0.68%  "
$ws.Range("D19").Formula = "'12.19"
$ws.Range("E19").Value = "  -2.35%  "
$ws.Range("D20").Formula = "'" + "0.0" + [char]0x2083 + "0904"
$ws.Range("E20").Value = "  -0.60%  "
$ws.Range("D21").Formula = "'5.92"
$ws.Range("E21").Value = "  -0.62%  "
$ws.Range("D22").Formula = "'67.17"
$ws.Range("E22").Value = "  -0.55%  "
$ws.Range("D23").Formula = "'234.91"
$ws.Range("E23").Value = "  -2.33%  "
$ws.Range("D24").Formula = "'2.57"
$ws.Range("E24").Value = "  -0.93%  "
$ws.Range("E25").Value = "  +0.47%  "
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("D27").Formula = "'23.38"
$ws.Range("E27").Value = "  -1.78%  "
$ws.Range("D28").Formula = "'36.66"
$ws.Range("E28").Value = "  -2.89%  "
$ws.Range("D29").Formula = "'2.12"
$ws.Range("E29").Value = "  +1.48%  "
$ws.Range("D30").Formula = "'9.49"
$ws.Range("E30").Value = "  +0.02%  "
$ws.Range("D31").Formula = "'164.62"
$ws.Range("E31").Value = "  +3.40%  "
$ws.Range("E32").Value = "  +0.03%  "
$ws.Range("D33").Formula = "'5.18"
$ws.Range("E33").Value = "  -0.65%  "
$ws.Range("E34").Value = "  -1.61%  "
$ws.Range("E35").Value = "  +3.06%  "
$ws.Range("D36").Formula = "'0.0722"
$ws.Range("E36").Value = "  -2.49%  "
$ws.Range("E37").Value = "  -0.60%  "
$ws.Range("E38").Value = "  -0.02%  "
$ws.Range("E39").Value = "  -2.45%  "
$ws.Range("D40").Formula = "'1.79"
$ws.Range("E40").Value = "  -2.83%  "
$ws.Range("D41").Formula = "'4.09"
$ws.Range("E41").Value = "  -0.23%  "
$ws.Range("D42").Formula = "'1.938.48"
$ws.Range("E42").Value = "  -3.00%  "
$ws.Range("B43").Value = "ApeXProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D43").Formula = "'2.22"
$ws.Range("E43").Value = "  -8.15%  "
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").Formula = "'0.0282"
$ws.Range("E44").Value = "  -1.32%  "
$ws.Range("E45").Value = "  -2.54%  "
$ws.Range("D46").Formula = "'2.91"
$ws.Range("E46").Value = "  -2.59%  "
$ws.Range("D47").Formula = "'9.65"
$ws.Range("E47").Value = "  -3.09%  "
$ws.Range("D48").Formula = "'53.71"
$ws.Range("E48").Value = "  +1.47%  "
$ws.Range("D49").Formula = "'2.473.11"
$ws.Range("E49").Value = "  -0.98%  "
$ws.Range("D50").Formula = "'71.23"
$ws.Range("E50").Value = "  -1.30%  "
$ws.Range("D51").Formula = "'91.25"
$ws.Range("E51").Value = "  -0.73%  "
